$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing row 2 data (becomes "Beşiktaş" -> "Kadıköy" with E2=0)
$ws.Range("B2").Value = "Beşiktaş"
$ws.Range("C2").Value = "Kadıköy"
$ws.Range("E2").Value = 0

# Insert new row 3 (original row 2 content: İstanbul Havalimanı -> Sabiha Gökçen)
$ws.Rows.Item(3).Insert()
$ws.Range("A3").Value = 2
$ws.Range("B3").Value = "İstanbul Havalimanı (IST)"
$ws.Range("C3").Value = "Sabiha Gökçen Havalimanı (SAW)"
$ws.Range("D3").Value = 0
$ws.Range("E3").Value = 12
$ws.Range("F3").Value = 0
$ws.Range("G3").Value = 0
$ws.Range("H3").Value = $true
$ws.Range("I3").Value = 0
$ws.Range("J3").Value = ""
$ws.Range("K3").Value = ""

# Insert new row 4 (İstanbul Havalimanı -> Sultanahmet, Fatih)
$ws.Rows.Item(4).Insert()
$ws.Range("A4").Value = 3
$ws.Range("B4").Value = "İstanbul Havalimanı (IST)"
$ws.Range("C4").Value = "Sultanahmet, Fatih"
$ws.Range("D4").Value = 0
$ws.Range("E4").Value = 2077
$ws.Range("F4").Value = 2290
$ws.Range("G4").Value = 2885
$ws.Range("H4").Value = $true
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 2127
$ws.Range("K4").Value = "Rakipten 50 TL ucuz (istanbulshuttleport.com)"
